$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new (blank) column before column N, shifting the old
# "Late"/"Outstanding" data columns (N:P) one place to the right (O:Q)
$ws.Columns("N:N").Insert()

# The new column N picks up the same width as column M, which is what
# Excel does automatically when inserting a column (format copied from
# the column to the left)
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Leave the selection where the author ended up after the edit
$ws.Range("S8").Select() | Out-Null
